$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.509.79'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').Value = '3.311.37'
$ws.Range('E3').Value = '  -0.29%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''587.12'
$ws.Range('E5').Value = '  +2.23%  '
$ws.Range('D6').Value = '''180.93'
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('D7').Value = '''0.654'
$ws.Range('E7').Value = '  +5.91%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '3.310.65'
$ws.Range('E9').Value = '  -0.27%  '
$ws.Range('E10').Value = '  -0.82%  '
$ws.Range('D11').Value = '''6.85'
$ws.Range('E11').Value = '  +2.76%  '
$ws.Range('E12').Value = '  +0.20%  '
$ws.Range('D13').Value = '3.891.82'
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('E14').Value = '  -2.70%  '
$ws.Range('D15').Value = '66.465.10'
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').Value = '''26.68'
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('D17').Value = '3.324.15'
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('D19').Value = '''424.02'
$ws.Range('E19').Value = '  -2.62%  '
$ws.Range('D20').Value = '''13.14'
$ws.Range('E20').Value = '  -3.16%  '
$ws.Range('E21').Value = '  -3.27%  '
$ws.Range('D22').Value = '''7.33'
$ws.Range('E22').Value = '  -3.16%  '
$ws.Range('D23').Value = '''71.53'
$ws.Range('E23').Value = '  -2.51%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').Value = '3.466.40'
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D26').Value = '''0.515'
$ws.Range('E26').Value = '  -1.12%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Value = '''0.205'
$ws.Range('E27').Value = '  +5.50%  '
$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').Value = '''0.0000115'
$ws.Range('E28').Value = '  -1.21%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = '''9.12'
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = '''0.999'
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '''1.92'
$ws.Range('E31').Value = '  -1.82%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = '''22.39'
$ws.Range('E32').Value = '  -1.51%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').Value = '''1.00'
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = '''5.18'
$ws.Range('E34').Value = '  -1.45%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').Value = '''6.59'
$ws.Range('E35').Value = '  -2.67%  '
$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D36').Value = '''1.19'
$ws.Range('E36').Value = '  -2.75%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '''160.15'
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '''1.43'
$ws.Range('E38').Value = '  -3.55%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '''1.80'
$ws.Range('E39').Value = '  +0.20%  '
$ws.Range('D40').Value = '2.858.36'
$ws.Range('E40').Value = '  +0.72%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '''26.35'
$ws.Range('E41').Value = '  -4.67%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = '''4.33'
$ws.Range('E42').Value = '  -2.59%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').Value = '''0.757'
$ws.Range('E43').Value = '  -4.22%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = '''39.71'
$ws.Range('E44').Value = '  -1.99%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').Value = '''0.0658'
$ws.Range('E45').Value = '  -1.37%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = '''5.89'
$ws.Range('E46').Value = '  -4.80%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').Value = '''2.30'
$ws.Range('E47').Value = '  -2.22%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '''23.12'
$ws.Range('E48').Value = '  -4.61%  '
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').Value = '''310.58'
$ws.Range('E49').Value = '  -4.60%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '''0.0273'
$ws.Range('E50').Value = '  +0.28%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '''0.104'
$ws.Range('E51').Value = '  +2.31%  '
